$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values
$ws.Range("C2").Value = 2
$ws.Range("C4").Value = 3
$ws.Range("A5").Value = "-"
$ws.Range("B5").Value = "-"
$ws.Range("C6").Value = 3

# Update the selection: selected range A2:B4 (active cell B4)
$ws.Range("A2:B4").Select()
